$wb = $excel.ActiveWorkbook

# ----- Sheet ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 3144.5454
$ws.Range("I80").Value = 2380.8
$ws.Range("J80").Value = 3781
$ws.Range("K80").Value = 7142.400000000001
$ws.Range("L80").Value = 11343
$ws.Range("M80").Value = -6144.400000000001
$ws.Range("N80").Value = -13339
$ws.Range("H83").Value = 3144.5454
$ws.Range("I83").Value = 2380.8
$ws.Range("J83").Value = 3781
$ws.Range("K83").Value = 21427.2
$ws.Range("L83").Value = 34029
$ws.Range("M83").Value = -16435.2
$ws.Range("N83").Value = -44013
$ws.Range("H86").Value = 6093.077
$ws.Range("I86").Value = 5458
$ws.Range("K86").Value = 5458
$ws.Range("M86").Value = -4335
$ws.Range("H89").Value = 6093.077
$ws.Range("I89").Value = 5458
$ws.Range("K89").Value = 27290
$ws.Range("M89").Value = -21674
$ws.Range("H112").Value = 1422
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1422
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 4266
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -6482
$ws.Range("H118").Value = 482.41666
$ws.Range("I118").Value = 332.1111
$ws.Range("J118").Value = 933.3333
$ws.Range("K118").Value = 996.3333
$ws.Range("L118").Value = 2799.9999
$ws.Range("M118").Value = 660.6667
$ws.Range("N118").Value = -6113.9999
$ws.Range("H125").Value = 972.8570999999999
$ws.Range("J125").Value = 972.8570999999999
$ws.Range("L125").Value = 8755.713899999999
$ws.Range("N125").Value = -13675.7139

# ----- Sheet ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4106.2344
$ws.Range("I32").Value = 2772.736
$ws.Range("J32").Value = 14774.223
$ws.Range("K32").Value = 2772.736
$ws.Range("L32").Value = 14774.223
$ws.Range("M32").Value = -2485.736
$ws.Range("N32").Value = -15348.223
$ws.Range("H61").Value = 1112.8636
$ws.Range("I61").Value = 1112.8636
$ws.Range("K61").Value = 1112.8636
$ws.Range("M61").Value = -900.8635999999999
$ws.Range("H74").Value = 48767.305
$ws.Range("I74").Value = 50886.7
$ws.Range("K74").Value = 50886.7
$ws.Range("M74").Value = -50012.7
$ws.Range("H77").Value = 48767.305
$ws.Range("I77").Value = 50886.7
$ws.Range("K77").Value = 254433.5
$ws.Range("M77").Value = -250065.5
$ws.Range("H88").Value = 2490
$ws.Range("I88").Value = 2320
$ws.Range("K88").Value = 2320
$ws.Range("M88").Value = -1914
$ws.Range("H91").Value = 2490
$ws.Range("I91").Value = 2320
$ws.Range("K91").Value = 2320
$ws.Range("M91").Value = -916
$ws.Range("H97").Value = 916.0769
$ws.Range("I97").Value = 634.44446
$ws.Range("J97").Value = 1549.75
$ws.Range("K97").Value = 634.44446
$ws.Range("L97").Value = 1549.75
$ws.Range("M97").Value = -138.44446
$ws.Range("N97").Value = -2541.75
$ws.Range("H122").Value = 1871.0769
$ws.Range("I122").Value = 1424.8889
$ws.Range("K122").Value = 4274.6667
$ws.Range("M122").Value = -1824.6667
$ws.Range("H132").Value = 2105.8103
$ws.Range("I132").Value = 1856.82
$ws.Range("J132").Value = 3662
$ws.Range("K132").Value = 5570.46
$ws.Range("L132").Value = 10986
$ws.Range("M132").Value = -3040.46
$ws.Range("N132").Value = -16046
$ws.Range("H136").Value = 1112.8636
$ws.Range("I136").Value = 1112.8636
$ws.Range("K136").Value = 3338.5908
$ws.Range("M136").Value = -788.5907999999999

# ----- Sheet BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1507.8485
$ws.Range("I107").Value = 1478.8695
$ws.Range("J107").Value = 1574.5
$ws.Range("K107").Value = 1478.8695
$ws.Range("L107").Value = 1574.5
$ws.Range("M107").Value = 441.1305
$ws.Range("N107").Value = -5414.5
$ws.Range("H134").Value = 3912.6765
$ws.Range("I134").Value = 2983.9565
$ws.Range("J134").Value = 5854.5454
$ws.Range("K134").Value = 8951.869499999999
$ws.Range("L134").Value = 17563.6362
$ws.Range("M134").Value = -6416.869499999999
$ws.Range("N134").Value = -22633.6362

# ----- Sheet CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21740704
$ws.Range("I31").Value = 29412702
$ws.Range("J31").Value = 3378.75
$ws.Range("K31").Value = 29412702
$ws.Range("L31").Value = 3378.75
$ws.Range("M31").Value = -29412407
$ws.Range("N31").Value = -3968.75
$ws.Range("H34").Value = 21740704
$ws.Range("I34").Value = 29412702
$ws.Range("J34").Value = 3378.75
$ws.Range("K34").Value = 29412702
$ws.Range("L34").Value = 3378.75
$ws.Range("M34").Value = -29412500
$ws.Range("N34").Value = -3782.75
$ws.Range("H58").Value = 956.8889
$ws.Range("I58").Value = 867.8214
$ws.Range("J58").Value = 1268.625
$ws.Range("K58").Value = 867.8214
$ws.Range("L58").Value = 1268.625
$ws.Range("M58").Value = -664.8214
$ws.Range("N58").Value = -1674.625
$ws.Range("H136").Value = 956.8889
$ws.Range("I136").Value = 867.8214
$ws.Range("J136").Value = 1268.625
$ws.Range("K136").Value = 2603.4642
$ws.Range("L136").Value = 3805.875
$ws.Range("M136").Value = -53.46420000000035
$ws.Range("N136").Value = -8905.875

# ----- Sheet CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1377894.5
$ws.Range("I113").Value = 1894424.5
$ws.Range("J113").Value = 481.16666
$ws.Range("K113").Value = 5683273.5
$ws.Range("L113").Value = 1443.49998
$ws.Range("M113").Value = -5681103.5
$ws.Range("N113").Value = -5783.499980000001
$ws.Range("H131").Value = 879.26
$ws.Range("I131").Value = 587.6667
$ws.Range("J131").Value = 908.0989
$ws.Range("K131").Value = 1763.0001
$ws.Range("L131").Value = 2724.2967
$ws.Range("M131").Value = 3276.9999
$ws.Range("N131").Value = -12804.2967
$ws.Range("H136").Value = 38467100
$ws.Range("I136").Value = 83335170
$ws.Range("J136").Value = 8757.143
$ws.Range("K136").Value = 250005510
$ws.Range("L136").Value = 26271.429
$ws.Range("M136").Value = -250000410
$ws.Range("N136").Value = -36471.429
$ws.Range("H137").Value = 15742111
$ws.Range("I137").Value = 2299.0908
$ws.Range("J137").Value = 21152672
$ws.Range("K137").Value = 6897.2724
$ws.Range("L137").Value = 63458016
$ws.Range("M137").Value = -1797.2724
$ws.Range("N137").Value = -63468216

# ----- Sheet GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2504.3928
$ws.Range("I80").Value = 2495.875
$ws.Range("J80").Value = 2555.5
$ws.Range("K80").Value = 2495.875
$ws.Range("L80").Value = 2555.5
$ws.Range("M80").Value = -1497.875
$ws.Range("N80").Value = -4551.5
$ws.Range("H83").Value = 2504.3928
$ws.Range("I83").Value = 2495.875
$ws.Range("J83").Value = 2555.5
$ws.Range("K83").Value = 12479.375
$ws.Range("L83").Value = 12777.5
$ws.Range("M83").Value = -7487.375
$ws.Range("N83").Value = -22761.5
$ws.Range("H102").Value = 1506.65
$ws.Range("I102").Value = 1237.5
$ws.Range("J102").Value = 2134.6667
$ws.Range("K102").Value = 1237.5
$ws.Range("L102").Value = 2134.6667
$ws.Range("M102").Value = 384.5
$ws.Range("N102").Value = -5378.6667
$ws.Range("H132").Value = 3684.4517
$ws.Range("I132").Value = 3774.318
$ws.Range("J132").Value = 3464.7778
$ws.Range("K132").Value = 11322.954
$ws.Range("L132").Value = 10394.3334
$ws.Range("M132").Value = -8792.954000000002
$ws.Range("N132").Value = -15454.3334

# ----- Sheet LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1168
$ws.Range("I7").Value = 1168
$ws.Range("K7").Value = 1168
$ws.Range("M7").Value = -1056
$ws.Range("H126").Value = 1168
$ws.Range("I126").Value = 1168
$ws.Range("K126").Value = 3504
$ws.Range("M126").Value = -1034
$ws.Range("H136").Value = 16668251
$ws.Range("I136").Value = 17545422
$ws.Range("K136").Value = 52636266
$ws.Range("M136").Value = -52633716

# ----- Sheet WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 24571.477
$ws.Range("J57").Value = 24571.477
$ws.Range("L57").Value = 24571.477
$ws.Range("N57").Value = -26079.477
$ws.Range("H81").Value = 955.8889
$ws.Range("I81").Value = 1133.5
$ws.Range("J81").Value = 600.6667
$ws.Range("K81").Value = 2267
$ws.Range("L81").Value = 1201.3334
$ws.Range("M81").Value = -1206
$ws.Range("N81").Value = -3323.3334
$ws.Range("H84").Value = 955.8889
$ws.Range("I84").Value = 1133.5
$ws.Range("J84").Value = 600.6667
$ws.Range("K84").Value = 11335
$ws.Range("L84").Value = 6006.666999999999
$ws.Range("M84").Value = -6031
$ws.Range("N84").Value = -16614.667
$ws.Range("H112").Value = 29999.357
$ws.Range("J112").Value = 29999.357
$ws.Range("L112").Value = 29999.357
$ws.Range("N112").Value = -32953.357
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
